$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header (H1), matching the formatting of the
# existing header row (copy format from G1, then set the text).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# New "Save" data values for the two existing rows.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
